# "contoh" -> "Contoh 1", split across three runs: "C" | "ontoh" | " 1"
# (matches the commit "update report & rename silidia")

$d = $word.ActiveDocument

# --- 1. Capitalize the leading "c" -> "C" ------------------------------
$rFirst = $d.Range(0, 1)
$rFirst.Text = "C"

# --- 2. Append " 1" after "...ontoh" ------------------------------------
$rEnd = $d.Range(6, 6)
$rEnd.InsertAfter(" 1")

# Paragraph text is now "Contoh 1". Word's run-merge pass would otherwise
# coalesce identically-formatted text back into a single run, so force
# the run boundaries at "C" | "ontoh" | " 1" by round-tripping a
# temporary bookmark over each piece (adding, then immediately deleting,
# a bookmark splits the enclosing run at its edges without leaving any
# residual formatting/markup behind).
function Split-RunAt($rangeObj, $bookmarkName) {
    $d.Bookmarks.Add($bookmarkName, $rangeObj) | Out-Null
    $d.Bookmarks($bookmarkName).Delete()
}

Split-RunAt $d.Range(0, 1) "__split_C__"
Split-RunAt $d.Range(1, 6) "__split_ontoh__"
